$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "NEW Postman Hospital 1 "
$ws.Range("G3").Value = "NEW Postman Hospital 2"
$ws.Range("G4").Value = "NEW Postman Hospital 3"
$ws.Range("G5").Value = "NEW Postman Hospital 4"
$ws.Range("G6").Value = "NEW Postman Hospital 5"
